$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Magias" glossary section (rows 53-69): the second column held placeholder
# duplicate-looking entries (Fogo/Fogo2/Fogo3, Terra2/Terra3, Sombra2/Sombra3,
# Luz2, Água2/Água3, Vento2/Vento3) which are replaced with the real
# inflected Portuguese glossary terms.
$ws.Range("C70").ClearContents()

$ws.Range("C53").Value = "Chama"
$ws.Range("C54").Value = "Chamae"
$ws.Range("C57").Value = "Terrae"
$ws.Range("C60").Value = "Sombrae"
$ws.Range("C63").Value = "Luze"
$ws.Range("C58").Value = "Terraus"
$ws.Range("C55").Value = "Chamaus"
$ws.Range("C61").Value = "Sombraus"
$ws.Range("C68").Value = "Ventoe"
$ws.Range("C69").Value = "Vetous"
$ws.Range("C65").Value = "Aguae"
$ws.Range("C66").Value = "Aguaus"

# The column got a bit narrower after the content change and lost its
# "best fit" auto-sizing (now a fixed, manually-set width).
$ws.Columns.Item(3).ColumnWidth = 26.6

# Reflect the new scroll position / selection the author left the sheet in.
$excel.ActiveWindow.ScrollRow = 28
[void]$ws.Range("C67").Select()
